$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark from the very first paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Mark the "50% mark" conversation (between the "50%" and "75%" health
#    milestones) as already-completed by striking through MILO's grey
#    dialogue lines (and the bold/italic storyline-transfer heading).
# ---------------------------------------------------------------------------
$strikeParagraphs = @(81, 83, 85, 87, 89, 91, 93, 95, 97, 99, 101, 102, 105)

foreach ($idx in $strikeParagraphs) {
    $p = $d.Paragraphs($idx)
    $p.Range.Font.StrikeThrough = 1
}

# ---------------------------------------------------------------------------
# 3. Re-add the "_GoBack" bookmark around "Fine. Your loss, I suppose."
#    (the paragraph that now marks the end of the struck-through section).
# ---------------------------------------------------------------------------
$fineRange = $d.Content
$fineRange.Find.Execute("Fine. Your loss, I suppose.") | Out-Null
$d.Bookmarks.Add("_GoBack", $fineRange)
